$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7795
$ws.Range("I43").Value = 7795
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 7795
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -7726
$ws.Range("N43").Value = $null

$ws.Range("H69").Value = 1013
$ws.Range("I69").Value = 1013
$ws.Range("K69").Value = 3039
$ws.Range("M69").Value = -2165

$ws.Range("H72").Value = 1013
$ws.Range("I72").Value = 1013
$ws.Range("K72").Value = 9117
$ws.Range("M72").Value = -4749

$ws.Range("H80").Value = 646
$ws.Range("I80").Value = 547.2857
$ws.Range("J80").Value = 1337
$ws.Range("K80").Value = 1641.8571
$ws.Range("L80").Value = 4011
$ws.Range("M80").Value = -643.8571000000002
$ws.Range("N80").Value = -6007

$ws.Range("H83").Value = 646
$ws.Range("I83").Value = 547.2857
$ws.Range("J83").Value = 1337
$ws.Range("K83").Value = 4925.571300000001
$ws.Range("L83").Value = 12033
$ws.Range("M83").Value = 66.42869999999948
$ws.Range("N83").Value = -22017

$ws.Range("H137").Value = 5403.5806
$ws.Range("J137").Value = 5116.9414
$ws.Range("L137").Value = 15350.8242
$ws.Range("N137").Value = -20450.8242

$ws.Range("H138").Value = 4108.5713
$ws.Range("I138").Value = 2915.25
$ws.Range("J138").Value = 4341.4146
$ws.Range("K138").Value = 8745.75
$ws.Range("L138").Value = 13024.2438
$ws.Range("M138").Value = -3605.75
$ws.Range("N138").Value = -23304.2438

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2084.8948
$ws.Range("I61").Value = 1859.8823
$ws.Range("K61").Value = 1859.8823
$ws.Range("M61").Value = -1647.8823

$ws.Range("H74").Value = 398.7143
$ws.Range("I74").Value = 381.83334
$ws.Range("K74").Value = 381.83334
$ws.Range("M74").Value = 492.16666

$ws.Range("H77").Value = 398.7143
$ws.Range("I77").Value = 381.83334
$ws.Range("K77").Value = 1909.1667
$ws.Range("M77").Value = 2458.8333

$ws.Range("H97").Value = 330.63635
$ws.Range("I97").Value = 330.63635
$ws.Range("K97").Value = 330.63635
$ws.Range("M97").Value = 165.36365

$ws.Range("H132").Value = 3779.9285
$ws.Range("I132").Value = 2490.375
$ws.Range("J132").Value = 5499.3335
$ws.Range("K132").Value = 7471.125
$ws.Range("L132").Value = 16498.0005
$ws.Range("M132").Value = -4941.125
$ws.Range("N132").Value = -21558.0005

$ws.Range("H136").Value = 2084.8948
$ws.Range("I136").Value = 1859.8823
$ws.Range("K136").Value = 5579.6469
$ws.Range("M136").Value = -3029.6469

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 100000
$ws.Range("J42").Value = 100000
$ws.Range("L42").Value = 100000
$ws.Range("N42").Value = -100656

$ws.Range("H134").Value = 1898.8889
$ws.Range("I134").Value = 1898.8889
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5696.6667
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3161.6667
$ws.Range("N134").Value = $null

$ws.Range("H138").Value = 125000
$ws.Range("J138").Value = 125000
$ws.Range("L138").Value = 125000
$ws.Range("N138").Value = -135280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 999
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 999
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 999
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = -1225

$ws.Range("H22").Value = 2049
$ws.Range("I22").Value = 2049
$ws.Range("K22").Value = 2049
$ws.Range("M22").Value = -1699

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = $null

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = $null

$ws.Range("H99").Value = 1577.4445
$ws.Range("I99").Value = 1685.4286
$ws.Range("K99").Value = 1685.4286
$ws.Range("M99").Value = -187.4286

$ws.Range("H126").Value = 1577.4445
$ws.Range("I126").Value = 1685.4286
$ws.Range("K126").Value = 5056.2858
$ws.Range("M126").Value = -2586.2858

$ws.Range("H134").Value = 7389.1
$ws.Range("I134").Value = 7765.778
$ws.Range("K134").Value = 23297.334
$ws.Range("M134").Value = -20762.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1851907
$ws.Range("I2").Value = 2645544.5
$ws.Range("K2").Value = 15873267
$ws.Range("M2").Value = -15873154

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1001.5
$ws.Range("I122").Value = 1002
$ws.Range("K122").Value = 3006
$ws.Range("M122").Value = -556

$ws.Range("H126").Value = 2055
$ws.Range("I126").Value = 1637.75
$ws.Range("J126").Value = 2333.1667
$ws.Range("K126").Value = 4913.25
$ws.Range("L126").Value = 6999.500100000001
$ws.Range("M126").Value = -2443.25
$ws.Range("N126").Value = -11939.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 1009
$ws.Range("J17").Value = 1009
$ws.Range("L17").Value = 1009
$ws.Range("N17").Value = -1349

$ws.Range("H22").Value = 2450
$ws.Range("I22").Value = 2450
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2450
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2155
$ws.Range("N22").Value = $null

$ws.Range("H27").Value = 2450
$ws.Range("I27").Value = 2450
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2450
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -2343
$ws.Range("N27").Value = $null

$ws.Range("H32").Value = 500
$ws.Range("I32").Value = 500
$ws.Range("K32").Value = 500
$ws.Range("M32").Value = -183

$ws.Range("H40").Value = 4940.2
$ws.Range("I40").Value = 4940.2
$ws.Range("K40").Value = 4940.2
$ws.Range("M40").Value = -4804.2

$ws.Range("H46").Value = 3825
$ws.Range("I46").Value = 3700
$ws.Range("J46").Value = 3887.5
$ws.Range("K46").Value = 3700
$ws.Range("L46").Value = 3887.5
$ws.Range("M46").Value = -3512
$ws.Range("N46").Value = -4263.5

$ws.Range("H55").Value = 466.66666

$ws.Range("H61").Value = 4886.75
$ws.Range("I61").Value = 4870.7144
$ws.Range("K61").Value = 4870.7144
$ws.Range("M61").Value = -4668.7144

$ws.Range("H113").Value = 4886.75
$ws.Range("I113").Value = 4870.7144
$ws.Range("K113").Value = 4870.7144
$ws.Range("M113").Value = -2700.7144

$ws.Range("H136").Value = 33767.93
$ws.Range("I136").Value = 17536.857
$ws.Range("K136").Value = 52610.571
$ws.Range("M136").Value = -50060.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 391
$ws.Range("I107").Value = 369.6
$ws.Range("K107").Value = 1108.8
$ws.Range("M107").Value = 811.1999999999998

$ws.Range("H136").Value = 6694
$ws.Range("I136").Value = 6797.615
$ws.Range("K136").Value = 20392.845
$ws.Range("M136").Value = -17842.845

$ws.Range("H137").Value = 44997.5
$ws.Range("J137").Value = 44997.5
$ws.Range("L137").Value = 44997.5
$ws.Range("N137").Value = -55197.5
